# Auto-generated COM-interop script implementing the PlayerPerformance_3529.xlsx edit
$wb = $excel.ActiveWorkbook

$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

# --- ODI Batting: MATCH_CARD_LINK -> MATCH_CODE (header + values) ---
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"
$wsBatting.Cells.Item(2, 4).NumberFormat = "@"
$wsBatting.Cells.Item(2, 4).Value = "2781"
$wsBatting.Cells.Item(3, 4).NumberFormat = "@"
$wsBatting.Cells.Item(3, 4).Value = "2815"
$wsBatting.Cells.Item(4, 4).NumberFormat = "@"
$wsBatting.Cells.Item(4, 4).Value = "2816"
$wsBatting.Cells.Item(5, 4).NumberFormat = "@"
$wsBatting.Cells.Item(5, 4).Value = "2817"
$wsBatting.Cells.Item(6, 4).NumberFormat = "@"
$wsBatting.Cells.Item(6, 4).Value = "2835"
$wsBatting.Cells.Item(7, 4).NumberFormat = "@"
$wsBatting.Cells.Item(7, 4).Value = "3197"
$wsBatting.Cells.Item(8, 4).NumberFormat = "@"
$wsBatting.Cells.Item(8, 4).Value = "3198"
$wsBatting.Cells.Item(9, 4).NumberFormat = "@"
$wsBatting.Cells.Item(9, 4).Value = "3201"
$wsBatting.Cells.Item(10, 4).NumberFormat = "@"
$wsBatting.Cells.Item(10, 4).Value = "3203"
$wsBatting.Cells.Item(11, 4).NumberFormat = "@"
$wsBatting.Cells.Item(11, 4).Value = "3223"
$wsBatting.Cells.Item(12, 4).NumberFormat = "@"
$wsBatting.Cells.Item(12, 4).Value = "3225"
$wsBatting.Cells.Item(13, 4).NumberFormat = "@"
$wsBatting.Cells.Item(13, 4).Value = "3228"
$wsBatting.Cells.Item(14, 4).NumberFormat = "@"
$wsBatting.Cells.Item(14, 4).Value = "3230"
$wsBatting.Cells.Item(15, 4).NumberFormat = "@"
$wsBatting.Cells.Item(15, 4).Value = "3251"
$wsBatting.Cells.Item(16, 4).NumberFormat = "@"
$wsBatting.Cells.Item(16, 4).Value = "3267"
$wsBatting.Cells.Item(17, 4).NumberFormat = "@"
$wsBatting.Cells.Item(17, 4).Value = "3274"
$wsBatting.Cells.Item(18, 4).NumberFormat = "@"
$wsBatting.Cells.Item(18, 4).Value = "3277"
$wsBatting.Cells.Item(19, 4).NumberFormat = "@"
$wsBatting.Cells.Item(19, 4).Value = "3282"
$wsBatting.Cells.Item(20, 4).NumberFormat = "@"
$wsBatting.Cells.Item(20, 4).Value = "3287"
$wsBatting.Cells.Item(21, 4).NumberFormat = "@"
$wsBatting.Cells.Item(21, 4).Value = "3288"
$wsBatting.Cells.Item(22, 4).NumberFormat = "@"
$wsBatting.Cells.Item(22, 4).Value = "3289"
$wsBatting.Cells.Item(23, 4).NumberFormat = "@"
$wsBatting.Cells.Item(23, 4).Value = "3291"
$wsBatting.Cells.Item(24, 4).NumberFormat = "@"
$wsBatting.Cells.Item(24, 4).Value = "3372"
$wsBatting.Cells.Item(25, 4).NumberFormat = "@"
$wsBatting.Cells.Item(25, 4).Value = "3374"
$wsBatting.Cells.Item(26, 4).NumberFormat = "@"
$wsBatting.Cells.Item(26, 4).Value = "3399"
$wsBatting.Cells.Item(27, 4).NumberFormat = "@"
$wsBatting.Cells.Item(27, 4).Value = "3488"
$wsBatting.Cells.Item(28, 4).NumberFormat = "@"
$wsBatting.Cells.Item(28, 4).Value = "3489"
$wsBatting.Cells.Item(29, 4).NumberFormat = "@"
$wsBatting.Cells.Item(29, 4).Value = "3491"
$wsBatting.Cells.Item(30, 4).NumberFormat = "@"
$wsBatting.Cells.Item(30, 4).Value = "3500"
$wsBatting.Cells.Item(31, 4).NumberFormat = "@"
$wsBatting.Cells.Item(31, 4).Value = "3506"
$wsBatting.Cells.Item(32, 4).NumberFormat = "@"
$wsBatting.Cells.Item(32, 4).Value = "3509"
$wsBatting.Cells.Item(33, 4).NumberFormat = "@"
$wsBatting.Cells.Item(33, 4).Value = "3514"
$wsBatting.Cells.Item(34, 4).NumberFormat = "@"
$wsBatting.Cells.Item(34, 4).Value = "3531"
$wsBatting.Cells.Item(35, 4).NumberFormat = "@"
$wsBatting.Cells.Item(35, 4).Value = "3532"
$wsBatting.Cells.Item(36, 4).NumberFormat = "@"
$wsBatting.Cells.Item(36, 4).Value = "3533"
$wsBatting.Cells.Item(37, 4).NumberFormat = "@"
$wsBatting.Cells.Item(37, 4).Value = "3535"
$wsBatting.Cells.Item(38, 4).NumberFormat = "@"
$wsBatting.Cells.Item(38, 4).Value = "3569"
$wsBatting.Cells.Item(39, 4).NumberFormat = "@"
$wsBatting.Cells.Item(39, 4).Value = "3571"
$wsBatting.Cells.Item(40, 4).NumberFormat = "@"
$wsBatting.Cells.Item(40, 4).Value = "3574"
$wsBatting.Cells.Item(41, 4).NumberFormat = "@"
$wsBatting.Cells.Item(41, 4).Value = "3658"
$wsBatting.Cells.Item(42, 4).NumberFormat = "@"
$wsBatting.Cells.Item(42, 4).Value = "3662"
$wsBatting.Cells.Item(43, 4).NumberFormat = "@"
$wsBatting.Cells.Item(43, 4).Value = "3666"
$wsBatting.Cells.Item(44, 4).NumberFormat = "@"
$wsBatting.Cells.Item(44, 4).Value = "3677"
$wsBatting.Cells.Item(45, 4).NumberFormat = "@"
$wsBatting.Cells.Item(45, 4).Value = "3679"
$wsBatting.Cells.Item(46, 4).NumberFormat = "@"
$wsBatting.Cells.Item(46, 4).Value = "3713"
$wsBatting.Cells.Item(47, 4).NumberFormat = "@"
$wsBatting.Cells.Item(47, 4).Value = "3715"
$wsBatting.Cells.Item(48, 4).NumberFormat = "@"
$wsBatting.Cells.Item(48, 4).Value = "3717"
$wsBatting.Cells.Item(49, 4).NumberFormat = "@"
$wsBatting.Cells.Item(49, 4).Value = "3751"
$wsBatting.Cells.Item(50, 4).NumberFormat = "@"
$wsBatting.Cells.Item(50, 4).Value = "3757"
$wsBatting.Cells.Item(51, 4).NumberFormat = "@"
$wsBatting.Cells.Item(51, 4).Value = "3770"
$wsBatting.Cells.Item(52, 4).NumberFormat = "@"
$wsBatting.Cells.Item(52, 4).Value = "3772"
$wsBatting.Cells.Item(53, 4).NumberFormat = "@"
$wsBatting.Cells.Item(53, 4).Value = "3776"
$wsBatting.Cells.Item(54, 4).NumberFormat = "@"
$wsBatting.Cells.Item(54, 4).Value = "3789"
$wsBatting.Cells.Item(55, 4).NumberFormat = "@"
$wsBatting.Cells.Item(55, 4).Value = "3792"
$wsBatting.Cells.Item(56, 4).NumberFormat = "@"
$wsBatting.Cells.Item(56, 4).Value = "3797"
$wsBatting.Cells.Item(57, 4).NumberFormat = "@"
$wsBatting.Cells.Item(57, 4).Value = "3798"
$wsBatting.Cells.Item(58, 4).NumberFormat = "@"
$wsBatting.Cells.Item(58, 4).Value = "3799"
$wsBatting.Cells.Item(59, 4).NumberFormat = "@"
$wsBatting.Cells.Item(59, 4).Value = "3801"
$wsBatting.Cells.Item(60, 4).NumberFormat = "@"
$wsBatting.Cells.Item(60, 4).Value = "3802"
$wsBatting.Cells.Item(61, 4).NumberFormat = "@"
$wsBatting.Cells.Item(61, 4).Value = "3803"
$wsBatting.Cells.Item(62, 4).NumberFormat = "@"
$wsBatting.Cells.Item(62, 4).Value = "3836"
$wsBatting.Cells.Item(63, 4).NumberFormat = "@"
$wsBatting.Cells.Item(63, 4).Value = "3837"
$wsBatting.Cells.Item(64, 4).NumberFormat = "@"
$wsBatting.Cells.Item(64, 4).Value = "3838"
$wsBatting.Cells.Item(65, 4).NumberFormat = "@"
$wsBatting.Cells.Item(65, 4).Value = "3858"
$wsBatting.Cells.Item(66, 4).NumberFormat = "@"
$wsBatting.Cells.Item(66, 4).Value = "3859"
$wsBatting.Cells.Item(67, 4).NumberFormat = "@"
$wsBatting.Cells.Item(67, 4).Value = "3861"
$wsBatting.Cells.Item(68, 4).NumberFormat = "@"
$wsBatting.Cells.Item(68, 4).Value = "3863"
$wsBatting.Cells.Item(69, 4).NumberFormat = "@"
$wsBatting.Cells.Item(69, 4).Value = "3879"
$wsBatting.Cells.Item(70, 4).NumberFormat = "@"
$wsBatting.Cells.Item(70, 4).Value = "3883"
$wsBatting.Cells.Item(71, 4).NumberFormat = "@"
$wsBatting.Cells.Item(71, 4).Value = "3925"
$wsBatting.Cells.Item(72, 4).NumberFormat = "@"
$wsBatting.Cells.Item(72, 4).Value = "3926"
$wsBatting.Cells.Item(73, 4).NumberFormat = "@"
$wsBatting.Cells.Item(73, 4).Value = "3928"
$wsBatting.Cells.Item(74, 4).NumberFormat = "@"
$wsBatting.Cells.Item(74, 4).Value = "3939"
$wsBatting.Cells.Item(75, 4).NumberFormat = "@"
$wsBatting.Cells.Item(75, 4).Value = "3943"
$wsBatting.Cells.Item(76, 4).NumberFormat = "@"
$wsBatting.Cells.Item(76, 4).Value = "3944"
$wsBatting.Cells.Item(77, 4).NumberFormat = "@"
$wsBatting.Cells.Item(77, 4).Value = "3972"
$wsBatting.Cells.Item(78, 4).NumberFormat = "@"
$wsBatting.Cells.Item(78, 4).Value = "3981"
$wsBatting.Cells.Item(79, 4).NumberFormat = "@"
$wsBatting.Cells.Item(79, 4).Value = "4017"
$wsBatting.Cells.Item(80, 4).NumberFormat = "@"
$wsBatting.Cells.Item(80, 4).Value = "4034"
$wsBatting.Cells.Item(81, 4).NumberFormat = "@"
$wsBatting.Cells.Item(81, 4).Value = "4304"
$wsBatting.Cells.Item(82, 4).NumberFormat = "@"
$wsBatting.Cells.Item(82, 4).Value = "4308"
$wsBatting.Cells.Item(83, 4).NumberFormat = "@"
$wsBatting.Cells.Item(83, 4).Value = "4319"
$wsBatting.Cells.Item(84, 4).NumberFormat = "@"
$wsBatting.Cells.Item(84, 4).Value = "4324"
$wsBatting.Cells.Item(85, 4).NumberFormat = "@"
$wsBatting.Cells.Item(85, 4).Value = "4334"
$wsBatting.Cells.Item(86, 4).NumberFormat = "@"
$wsBatting.Cells.Item(86, 4).Value = "4337"
$wsBatting.Cells.Item(87, 4).NumberFormat = "@"
$wsBatting.Cells.Item(87, 4).Value = "4340"
$wsBatting.Cells.Item(88, 4).NumberFormat = "@"
$wsBatting.Cells.Item(88, 4).Value = "4349"
$wsBatting.Cells.Item(89, 4).NumberFormat = "@"
$wsBatting.Cells.Item(89, 4).Value = "4375"
$wsBatting.Cells.Item(90, 4).NumberFormat = "@"
$wsBatting.Cells.Item(90, 4).Value = "4376"
$wsBatting.Cells.Item(91, 4).NumberFormat = "@"
$wsBatting.Cells.Item(91, 4).Value = "4432"
$wsBatting.Cells.Item(92, 4).NumberFormat = "@"
$wsBatting.Cells.Item(92, 4).Value = "4434"

# --- ODI Batting: drop stray empty INNING_NUMBER cells ---
$wsBatting.Cells.Item(2, 2).Value = $null
$wsBatting.Cells.Item(11, 2).Value = $null
$wsBatting.Cells.Item(16, 2).Value = $null
$wsBatting.Cells.Item(17, 2).Value = $null
$wsBatting.Cells.Item(18, 2).Value = $null
$wsBatting.Cells.Item(20, 2).Value = $null
$wsBatting.Cells.Item(21, 2).Value = $null
$wsBatting.Cells.Item(24, 2).Value = $null
$wsBatting.Cells.Item(26, 2).Value = $null
$wsBatting.Cells.Item(37, 2).Value = $null
$wsBatting.Cells.Item(41, 2).Value = $null
$wsBatting.Cells.Item(48, 2).Value = $null
$wsBatting.Cells.Item(54, 2).Value = $null
$wsBatting.Cells.Item(59, 2).Value = $null
$wsBatting.Cells.Item(60, 2).Value = $null
$wsBatting.Cells.Item(62, 2).Value = $null
$wsBatting.Cells.Item(64, 2).Value = $null
$wsBatting.Cells.Item(65, 2).Value = $null
$wsBatting.Cells.Item(71, 2).Value = $null
$wsBatting.Cells.Item(75, 2).Value = $null
$wsBatting.Cells.Item(76, 2).Value = $null
$wsBatting.Cells.Item(79, 2).Value = $null
$wsBatting.Cells.Item(84, 2).Value = $null
$wsBatting.Cells.Item(86, 2).Value = $null

# --- ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE (header + values) ---
$wsBowling.Cells.Item(1, 2).Value = "MATCH_CODE"
$wsBowling.Cells.Item(2, 2).NumberFormat = "@"
$wsBowling.Cells.Item(2, 2).Value = "2781"
$wsBowling.Cells.Item(3, 2).NumberFormat = "@"
$wsBowling.Cells.Item(3, 2).Value = "2815"
$wsBowling.Cells.Item(4, 2).NumberFormat = "@"
$wsBowling.Cells.Item(4, 2).Value = "2816"
$wsBowling.Cells.Item(5, 2).NumberFormat = "@"
$wsBowling.Cells.Item(5, 2).Value = "2817"
$wsBowling.Cells.Item(6, 2).NumberFormat = "@"
$wsBowling.Cells.Item(6, 2).Value = "2835"
$wsBowling.Cells.Item(7, 2).NumberFormat = "@"
$wsBowling.Cells.Item(7, 2).Value = "3197"
$wsBowling.Cells.Item(8, 2).NumberFormat = "@"
$wsBowling.Cells.Item(8, 2).Value = "3198"
$wsBowling.Cells.Item(9, 2).NumberFormat = "@"
$wsBowling.Cells.Item(9, 2).Value = "3201"
$wsBowling.Cells.Item(10, 2).NumberFormat = "@"
$wsBowling.Cells.Item(10, 2).Value = "3203"
$wsBowling.Cells.Item(11, 2).NumberFormat = "@"
$wsBowling.Cells.Item(11, 2).Value = "3225"
$wsBowling.Cells.Item(12, 2).NumberFormat = "@"
$wsBowling.Cells.Item(12, 2).Value = "3228"
$wsBowling.Cells.Item(13, 2).NumberFormat = "@"
$wsBowling.Cells.Item(13, 2).Value = "3230"
$wsBowling.Cells.Item(14, 2).NumberFormat = "@"
$wsBowling.Cells.Item(14, 2).Value = "3251"
$wsBowling.Cells.Item(15, 2).NumberFormat = "@"
$wsBowling.Cells.Item(15, 2).Value = "3267"
$wsBowling.Cells.Item(16, 2).NumberFormat = "@"
$wsBowling.Cells.Item(16, 2).Value = "3274"
$wsBowling.Cells.Item(17, 2).NumberFormat = "@"
$wsBowling.Cells.Item(17, 2).Value = "3277"
$wsBowling.Cells.Item(18, 2).NumberFormat = "@"
$wsBowling.Cells.Item(18, 2).Value = "3282"
$wsBowling.Cells.Item(19, 2).NumberFormat = "@"
$wsBowling.Cells.Item(19, 2).Value = "3287"
$wsBowling.Cells.Item(20, 2).NumberFormat = "@"
$wsBowling.Cells.Item(20, 2).Value = "3288"
$wsBowling.Cells.Item(21, 2).NumberFormat = "@"
$wsBowling.Cells.Item(21, 2).Value = "3289"
$wsBowling.Cells.Item(22, 2).NumberFormat = "@"
$wsBowling.Cells.Item(22, 2).Value = "3291"
$wsBowling.Cells.Item(23, 2).NumberFormat = "@"
$wsBowling.Cells.Item(23, 2).Value = "3372"
$wsBowling.Cells.Item(24, 2).NumberFormat = "@"
$wsBowling.Cells.Item(24, 2).Value = "3374"
$wsBowling.Cells.Item(25, 2).NumberFormat = "@"
$wsBowling.Cells.Item(25, 2).Value = "3399"
$wsBowling.Cells.Item(26, 2).NumberFormat = "@"
$wsBowling.Cells.Item(26, 2).Value = "3488"
$wsBowling.Cells.Item(27, 2).NumberFormat = "@"
$wsBowling.Cells.Item(27, 2).Value = "3489"
$wsBowling.Cells.Item(28, 2).NumberFormat = "@"
$wsBowling.Cells.Item(28, 2).Value = "3491"
$wsBowling.Cells.Item(29, 2).NumberFormat = "@"
$wsBowling.Cells.Item(29, 2).Value = "3500"
$wsBowling.Cells.Item(30, 2).NumberFormat = "@"
$wsBowling.Cells.Item(30, 2).Value = "3506"
$wsBowling.Cells.Item(31, 2).NumberFormat = "@"
$wsBowling.Cells.Item(31, 2).Value = "3509"
$wsBowling.Cells.Item(32, 2).NumberFormat = "@"
$wsBowling.Cells.Item(32, 2).Value = "3514"
$wsBowling.Cells.Item(33, 2).NumberFormat = "@"
$wsBowling.Cells.Item(33, 2).Value = "3531"
$wsBowling.Cells.Item(34, 2).NumberFormat = "@"
$wsBowling.Cells.Item(34, 2).Value = "3532"
$wsBowling.Cells.Item(35, 2).NumberFormat = "@"
$wsBowling.Cells.Item(35, 2).Value = "3533"
$wsBowling.Cells.Item(36, 2).NumberFormat = "@"
$wsBowling.Cells.Item(36, 2).Value = "3535"
$wsBowling.Cells.Item(37, 2).NumberFormat = "@"
$wsBowling.Cells.Item(37, 2).Value = "3569"
$wsBowling.Cells.Item(38, 2).NumberFormat = "@"
$wsBowling.Cells.Item(38, 2).Value = "3571"
$wsBowling.Cells.Item(39, 2).NumberFormat = "@"
$wsBowling.Cells.Item(39, 2).Value = "3574"
$wsBowling.Cells.Item(40, 2).NumberFormat = "@"
$wsBowling.Cells.Item(40, 2).Value = "3658"
$wsBowling.Cells.Item(41, 2).NumberFormat = "@"
$wsBowling.Cells.Item(41, 2).Value = "3662"
$wsBowling.Cells.Item(42, 2).NumberFormat = "@"
$wsBowling.Cells.Item(42, 2).Value = "3666"
$wsBowling.Cells.Item(43, 2).NumberFormat = "@"
$wsBowling.Cells.Item(43, 2).Value = "3677"
$wsBowling.Cells.Item(44, 2).NumberFormat = "@"
$wsBowling.Cells.Item(44, 2).Value = "3679"
$wsBowling.Cells.Item(45, 2).NumberFormat = "@"
$wsBowling.Cells.Item(45, 2).Value = "3713"
$wsBowling.Cells.Item(46, 2).NumberFormat = "@"
$wsBowling.Cells.Item(46, 2).Value = "3715"
$wsBowling.Cells.Item(47, 2).NumberFormat = "@"
$wsBowling.Cells.Item(47, 2).Value = "3717"
$wsBowling.Cells.Item(48, 2).NumberFormat = "@"
$wsBowling.Cells.Item(48, 2).Value = "3751"
$wsBowling.Cells.Item(49, 2).NumberFormat = "@"
$wsBowling.Cells.Item(49, 2).Value = "3757"
$wsBowling.Cells.Item(50, 2).NumberFormat = "@"
$wsBowling.Cells.Item(50, 2).Value = "3770"
$wsBowling.Cells.Item(51, 2).NumberFormat = "@"
$wsBowling.Cells.Item(51, 2).Value = "3772"
$wsBowling.Cells.Item(52, 2).NumberFormat = "@"
$wsBowling.Cells.Item(52, 2).Value = "3776"
$wsBowling.Cells.Item(53, 2).NumberFormat = "@"
$wsBowling.Cells.Item(53, 2).Value = "3789"
$wsBowling.Cells.Item(54, 2).NumberFormat = "@"
$wsBowling.Cells.Item(54, 2).Value = "3792"
$wsBowling.Cells.Item(55, 2).NumberFormat = "@"
$wsBowling.Cells.Item(55, 2).Value = "3797"
$wsBowling.Cells.Item(56, 2).NumberFormat = "@"
$wsBowling.Cells.Item(56, 2).Value = "3798"
$wsBowling.Cells.Item(57, 2).NumberFormat = "@"
$wsBowling.Cells.Item(57, 2).Value = "3799"
$wsBowling.Cells.Item(58, 2).NumberFormat = "@"
$wsBowling.Cells.Item(58, 2).Value = "3801"
$wsBowling.Cells.Item(59, 2).NumberFormat = "@"
$wsBowling.Cells.Item(59, 2).Value = "3802"
$wsBowling.Cells.Item(60, 2).NumberFormat = "@"
$wsBowling.Cells.Item(60, 2).Value = "3803"
$wsBowling.Cells.Item(61, 2).NumberFormat = "@"
$wsBowling.Cells.Item(61, 2).Value = "3836"
$wsBowling.Cells.Item(62, 2).NumberFormat = "@"
$wsBowling.Cells.Item(62, 2).Value = "3837"
$wsBowling.Cells.Item(63, 2).NumberFormat = "@"
$wsBowling.Cells.Item(63, 2).Value = "3838"
$wsBowling.Cells.Item(64, 2).NumberFormat = "@"
$wsBowling.Cells.Item(64, 2).Value = "3858"
$wsBowling.Cells.Item(65, 2).NumberFormat = "@"
$wsBowling.Cells.Item(65, 2).Value = "3859"
$wsBowling.Cells.Item(66, 2).NumberFormat = "@"
$wsBowling.Cells.Item(66, 2).Value = "3861"
$wsBowling.Cells.Item(67, 2).NumberFormat = "@"
$wsBowling.Cells.Item(67, 2).Value = "3863"
$wsBowling.Cells.Item(68, 2).NumberFormat = "@"
$wsBowling.Cells.Item(68, 2).Value = "3879"
$wsBowling.Cells.Item(69, 2).NumberFormat = "@"
$wsBowling.Cells.Item(69, 2).Value = "3883"
$wsBowling.Cells.Item(70, 2).NumberFormat = "@"
$wsBowling.Cells.Item(70, 2).Value = "3925"
$wsBowling.Cells.Item(71, 2).NumberFormat = "@"
$wsBowling.Cells.Item(71, 2).Value = "3926"
$wsBowling.Cells.Item(72, 2).NumberFormat = "@"
$wsBowling.Cells.Item(72, 2).Value = "3928"
$wsBowling.Cells.Item(73, 2).NumberFormat = "@"
$wsBowling.Cells.Item(73, 2).Value = "3939"
$wsBowling.Cells.Item(74, 2).NumberFormat = "@"
$wsBowling.Cells.Item(74, 2).Value = "3943"
$wsBowling.Cells.Item(75, 2).NumberFormat = "@"
$wsBowling.Cells.Item(75, 2).Value = "3944"
$wsBowling.Cells.Item(76, 2).NumberFormat = "@"
$wsBowling.Cells.Item(76, 2).Value = "3972"
$wsBowling.Cells.Item(77, 2).NumberFormat = "@"
$wsBowling.Cells.Item(77, 2).Value = "3981"
$wsBowling.Cells.Item(78, 2).NumberFormat = "@"
$wsBowling.Cells.Item(78, 2).Value = "4017"
$wsBowling.Cells.Item(79, 2).NumberFormat = "@"
$wsBowling.Cells.Item(79, 2).Value = "4034"
$wsBowling.Cells.Item(80, 2).NumberFormat = "@"
$wsBowling.Cells.Item(80, 2).Value = "4304"
$wsBowling.Cells.Item(81, 2).NumberFormat = "@"
$wsBowling.Cells.Item(81, 2).Value = "4308"
$wsBowling.Cells.Item(82, 2).NumberFormat = "@"
$wsBowling.Cells.Item(82, 2).Value = "4319"
$wsBowling.Cells.Item(83, 2).NumberFormat = "@"
$wsBowling.Cells.Item(83, 2).Value = "4324"
$wsBowling.Cells.Item(84, 2).NumberFormat = "@"
$wsBowling.Cells.Item(84, 2).Value = "4334"
$wsBowling.Cells.Item(85, 2).NumberFormat = "@"
$wsBowling.Cells.Item(85, 2).Value = "4337"
$wsBowling.Cells.Item(86, 2).NumberFormat = "@"
$wsBowling.Cells.Item(86, 2).Value = "4340"
$wsBowling.Cells.Item(87, 2).NumberFormat = "@"
$wsBowling.Cells.Item(87, 2).Value = "4349"
$wsBowling.Cells.Item(88, 2).NumberFormat = "@"
$wsBowling.Cells.Item(88, 2).Value = "4375"
$wsBowling.Cells.Item(89, 2).NumberFormat = "@"
$wsBowling.Cells.Item(89, 2).Value = "4376"
$wsBowling.Cells.Item(90, 2).NumberFormat = "@"
$wsBowling.Cells.Item(90, 2).Value = "4432"
$wsBowling.Cells.Item(91, 2).NumberFormat = "@"
$wsBowling.Cells.Item(91, 2).Value = "4434"

# --- Add 'Player Info' sheet (inserted before 'ODI Batting') ---
$wsPlayerInfo = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$wsPlayerInfo.Name = "Player Info"
$wsPlayerInfo.Cells.Item(1, 1).Value = "ID"
$wsPlayerInfo.Cells.Item(1, 2).Value = "NAME"
$wsPlayerInfo.Cells.Item(1, 3).Value = "BATTING_HAND"
$wsPlayerInfo.Cells.Item(1, 4).Value = "BOWL_STYLE"
$wsPlayerInfo.Cells.Item(2, 1).NumberFormat = "@"
$wsPlayerInfo.Cells.Item(2, 1).Value = "3529"
$wsPlayerInfo.Cells.Item(2, 2).NumberFormat = "@"
$wsPlayerInfo.Cells.Item(2, 2).Value = "Wahab Riaz"
$wsPlayerInfo.Cells.Item(2, 3).NumberFormat = "@"
$wsPlayerInfo.Cells.Item(2, 3).Value = "Right Handed"
$wsPlayerInfo.Cells.Item(2, 4).NumberFormat = "@"
$wsPlayerInfo.Cells.Item(2, 4).Value = "Left Arm Fast"

# --- Add 'ODI Batting Extra' sheet (inserted after 'ODI Bowling') ---
$wsBattingExtra = $wb.Worksheets.Add($null, $wb.Worksheets.Item("ODI Bowling"))
$wsBattingExtra.Name = "ODI Batting Extra"
$wsBattingExtra.Cells.Item(1, 1).Value = "MATCH_CODE"
$wsBattingExtra.Cells.Item(1, 2).Value = "BATTING_POSITION"
$wsBattingExtra.Cells.Item(1, 3).Value = "NUM_4"
$wsBattingExtra.Cells.Item(1, 4).Value = "NUM_6"
$wsBattingExtra.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$wsBattingExtra.Cells.Item(1, 6).Value = "MAN_OF_MATCH"
$wsBattingExtra.Cells.Item(2, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(2, 1).Value = "3928"
$wsBattingExtra.Cells.Item(2, 2).Value = 9
$wsBattingExtra.Cells.Item(2, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(2, 3).Value = "0"
$wsBattingExtra.Cells.Item(2, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(2, 4).Value = "0"
$wsBattingExtra.Cells.Item(2, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(2, 5).Value = "5.09%"
$wsBattingExtra.Cells.Item(2, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(2, 6).Value = "NO"
$wsBattingExtra.Cells.Item(3, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(3, 1).Value = "3939"
$wsBattingExtra.Cells.Item(3, 2).Value = 9
$wsBattingExtra.Cells.Item(3, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(3, 3).Value = "0"
$wsBattingExtra.Cells.Item(3, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(3, 4).Value = "0"
$wsBattingExtra.Cells.Item(3, 5).Value = $null
$wsBattingExtra.Cells.Item(3, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(3, 6).Value = "NO"
$wsBattingExtra.Cells.Item(4, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(4, 1).Value = "3943"
$wsBattingExtra.Cells.Item(4, 2).Value = 9
$wsBattingExtra.Cells.Item(4, 3).Value = $null
$wsBattingExtra.Cells.Item(4, 4).Value = $null
$wsBattingExtra.Cells.Item(4, 5).Value = $null
$wsBattingExtra.Cells.Item(4, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(4, 6).Value = "NO"
$wsBattingExtra.Cells.Item(5, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(5, 1).Value = "3944"
$wsBattingExtra.Cells.Item(5, 2).Value = 10
$wsBattingExtra.Cells.Item(5, 3).Value = $null
$wsBattingExtra.Cells.Item(5, 4).Value = $null
$wsBattingExtra.Cells.Item(5, 5).Value = $null
$wsBattingExtra.Cells.Item(5, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(5, 6).Value = "NO"
$wsBattingExtra.Cells.Item(6, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(6, 1).Value = "3972"
$wsBattingExtra.Cells.Item(6, 2).Value = 10
$wsBattingExtra.Cells.Item(6, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(6, 3).Value = "0"
$wsBattingExtra.Cells.Item(6, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(6, 4).Value = "0"
$wsBattingExtra.Cells.Item(6, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(6, 5).Value = "3.41%"
$wsBattingExtra.Cells.Item(6, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(6, 6).Value = "NO"
$wsBattingExtra.Cells.Item(7, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(7, 1).Value = "3981"
$wsBattingExtra.Cells.Item(7, 2).Value = 9
$wsBattingExtra.Cells.Item(7, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(7, 3).Value = "1"
$wsBattingExtra.Cells.Item(7, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(7, 4).Value = "0"
$wsBattingExtra.Cells.Item(7, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(7, 5).Value = "5.45%"
$wsBattingExtra.Cells.Item(7, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(7, 6).Value = "NO"
$wsBattingExtra.Cells.Item(8, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(8, 1).Value = "4017"
$wsBattingExtra.Cells.Item(8, 2).Value = 11
$wsBattingExtra.Cells.Item(8, 3).Value = $null
$wsBattingExtra.Cells.Item(8, 4).Value = $null
$wsBattingExtra.Cells.Item(8, 5).Value = $null
$wsBattingExtra.Cells.Item(8, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(8, 6).Value = "NO"
$wsBattingExtra.Cells.Item(9, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(9, 1).Value = "4034"
$wsBattingExtra.Cells.Item(9, 2).Value = 10
$wsBattingExtra.Cells.Item(9, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(9, 3).Value = "1"
$wsBattingExtra.Cells.Item(9, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(9, 4).Value = "2"
$wsBattingExtra.Cells.Item(9, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(9, 5).Value = "17.14%"
$wsBattingExtra.Cells.Item(9, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(9, 6).Value = "NO"
$wsBattingExtra.Cells.Item(10, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(10, 1).Value = "4304"
$wsBattingExtra.Cells.Item(10, 2).Value = $null
$wsBattingExtra.Cells.Item(10, 3).Value = $null
$wsBattingExtra.Cells.Item(10, 4).Value = $null
$wsBattingExtra.Cells.Item(10, 5).Value = $null
$wsBattingExtra.Cells.Item(10, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(10, 6).Value = "NO"
$wsBattingExtra.Cells.Item(11, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(11, 1).Value = "4308"
$wsBattingExtra.Cells.Item(11, 2).Value = 9
$wsBattingExtra.Cells.Item(11, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(11, 3).Value = "2"
$wsBattingExtra.Cells.Item(11, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(11, 4).Value = "3"
$wsBattingExtra.Cells.Item(11, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(11, 5).Value = "16.92%"
$wsBattingExtra.Cells.Item(11, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(11, 6).Value = "NO"
$wsBattingExtra.Cells.Item(12, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(12, 1).Value = "4319"
$wsBattingExtra.Cells.Item(12, 2).Value = 10
$wsBattingExtra.Cells.Item(12, 3).Value = $null
$wsBattingExtra.Cells.Item(12, 4).Value = $null
$wsBattingExtra.Cells.Item(12, 5).Value = $null
$wsBattingExtra.Cells.Item(12, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(12, 6).Value = "NO"
$wsBattingExtra.Cells.Item(13, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(13, 1).Value = "4324"
$wsBattingExtra.Cells.Item(13, 2).Value = 7
$wsBattingExtra.Cells.Item(13, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(13, 3).Value = "0"
$wsBattingExtra.Cells.Item(13, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(13, 4).Value = "0"
$wsBattingExtra.Cells.Item(13, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(13, 5).Value = "1.30%"
$wsBattingExtra.Cells.Item(13, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(13, 6).Value = "NO"
$wsBattingExtra.Cells.Item(14, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(14, 1).Value = "4334"
$wsBattingExtra.Cells.Item(14, 2).Value = $null
$wsBattingExtra.Cells.Item(14, 3).Value = $null
$wsBattingExtra.Cells.Item(14, 4).Value = $null
$wsBattingExtra.Cells.Item(14, 5).Value = $null
$wsBattingExtra.Cells.Item(14, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(14, 6).Value = "NO"
$wsBattingExtra.Cells.Item(15, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(15, 1).Value = "4337"
$wsBattingExtra.Cells.Item(15, 2).Value = 9
$wsBattingExtra.Cells.Item(15, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(15, 3).Value = "1"
$wsBattingExtra.Cells.Item(15, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(15, 4).Value = "1"
$wsBattingExtra.Cells.Item(15, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(15, 5).Value = "6.52%"
$wsBattingExtra.Cells.Item(15, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(15, 6).Value = "NO"
$wsBattingExtra.Cells.Item(16, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(16, 1).Value = "4340"
$wsBattingExtra.Cells.Item(16, 2).Value = 8
$wsBattingExtra.Cells.Item(16, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(16, 3).Value = "0"
$wsBattingExtra.Cells.Item(16, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(16, 4).Value = "0"
$wsBattingExtra.Cells.Item(16, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(16, 5).Value = "0.63%"
$wsBattingExtra.Cells.Item(16, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(16, 6).Value = "NO"
$wsBattingExtra.Cells.Item(17, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(17, 1).Value = "4349"
$wsBattingExtra.Cells.Item(17, 2).Value = 8
$wsBattingExtra.Cells.Item(17, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(17, 3).Value = "0"
$wsBattingExtra.Cells.Item(17, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(17, 4).Value = "0"
$wsBattingExtra.Cells.Item(17, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(17, 5).Value = "0.66%"
$wsBattingExtra.Cells.Item(17, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(17, 6).Value = "NO"
$wsBattingExtra.Cells.Item(18, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(18, 1).Value = "4375"
$wsBattingExtra.Cells.Item(18, 2).Value = 7
$wsBattingExtra.Cells.Item(18, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(18, 3).Value = "0"
$wsBattingExtra.Cells.Item(18, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(18, 4).Value = "0"
$wsBattingExtra.Cells.Item(18, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(18, 5).Value = "0.33%"
$wsBattingExtra.Cells.Item(18, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(18, 6).Value = "NO"
$wsBattingExtra.Cells.Item(19, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(19, 1).Value = "4376"
$wsBattingExtra.Cells.Item(19, 2).Value = $null
$wsBattingExtra.Cells.Item(19, 3).Value = $null
$wsBattingExtra.Cells.Item(19, 4).Value = $null
$wsBattingExtra.Cells.Item(19, 5).Value = $null
$wsBattingExtra.Cells.Item(19, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(19, 6).Value = "NO"
$wsBattingExtra.Cells.Item(20, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(20, 1).Value = "4432"
$wsBattingExtra.Cells.Item(20, 2).Value = 8
$wsBattingExtra.Cells.Item(20, 3).NumberFormat = "@"
$wsBattingExtra.Cells.Item(20, 3).Value = "3"
$wsBattingExtra.Cells.Item(20, 4).NumberFormat = "@"
$wsBattingExtra.Cells.Item(20, 4).Value = "3"
$wsBattingExtra.Cells.Item(20, 5).NumberFormat = "@"
$wsBattingExtra.Cells.Item(20, 5).Value = "18.71%"
$wsBattingExtra.Cells.Item(20, 6).NumberFormat = "@"
$wsBattingExtra.Cells.Item(20, 6).Value = "NO"
$wsBattingExtra.Cells.Item(21, 1).NumberFormat = "@"
$wsBattingExtra.Cells.Item(21, 1).Value = "4434"
$wsBattingExtra.Cells.Item(21, 2).Value = $null
$wsBattingExtra.Cells.Item(21, 3).Value = $null
$wsBattingExtra.Cells.Item(21, 4).Value = $null
$wsBattingExtra.Cells.Item(21, 5).Value = $null
$wsBattingExtra.Cells.Item(21, 6).Value = $null

